$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Allow spaces at the end of hearing time
$ws.Range("C2").Value = "10am           "
$ws.Range("C3").Value = "1:30pm                "

# Reflect the selected cell left in the sheet when the author saved
$ws.Range("C6").Select()
